{"js": "// The author's edit inserts a protocol-number prefix \"\u03a6.11.2/\" in front of\n// the \"${regionaldirect_protocol}\" merge field in the \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: ...\"\n// line of the document header table. (The surrounding w:lang=\"en-US\" /\n// run-splitting churn visible in the raw OOXML diff is a cosmetic\n// side-effect of Word's editor and has no effect on the rendered text.)\n\nconst body = context.document.body;\n\nconst results = body.search(\".: ${regionaldirect_protocol}\", { matchCase: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: ${regionaldirect_protocol}\" text to update.');\n}\n\n// Replace the matched range \".: ${regionaldirect_protocol}\" with\n// \".: \u03a6.11.2/${regionaldirect_protocol}\" so the full line reads\n// \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: \u03a6.11.2/${regionaldirect_protocol}\".\nresults.items[0].insertText(\".: \u03a6.11.2/${regionaldirect_protocol}\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The author's edit inserts a protocol-number prefix \"\u03a6.11.2/\" in front of\n# the \"${regionaldirect_protocol}\" merge field in the \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: ...\"\n# line of the document header table. (The surrounding w:lang=\"en-US\" /\n# run-splitting churn visible in the raw OOXML diff is a cosmetic\n# side-effect of Word's editor and has no effect on the rendered text.)\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.MatchWildcards = $false\n$found = $range.Find.Execute(\".: `${regionaldirect_protocol}\")\n\nif ($found) {\n    $range.Text = \".: \u03a6.11.2/`${regionaldirect_protocol}\"\n} else {\n    throw 'Could not find \"\u0391\u03c1. \u03a0\u03c1\u03c9\u03c4.: ${regionaldirect_protocol}\" text to update.'\n}\n"}
